# Auto-generated script applying 2024-10-18 YTD violent crime data updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6113
$ws.Range("I3").Value = 6072
$ws.Range("K3").Value = 6294
$ws.Range("C4").Value = 1539
$ws.Range("G4").Value = 1224
$ws.Range("J4").Value = 1499
$ws.Range("K4").Value = 1314
$ws.Range("K5").Value = 447
$ws.Range("K6").Value = 6940
$ws.Range("C7").Value = 22626
$ws.Range("G7").Value = 20028
$ws.Range("I7").Value = 20782
$ws.Range("J7").Value = 23205
$ws.Range("K7").Value = 21108

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("K5").Value = 11
$ws.Range("K6").Value = 21

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 380
$ws.Range("K3").Value = 423
$ws.Range("K5").Value = 42
$ws.Range("K6").Value = 469
$ws.Range("K7").Value = 1389

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 165
$ws.Range("K6").Value = 108

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 238
$ws.Range("K6").Value = 281
$ws.Range("K7").Value = 920

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 123
$ws.Range("K6").Value = 78

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K4").Value = 35
$ws.Range("K7").Value = 717

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 167
$ws.Range("K3").Value = 124
$ws.Range("K6").Value = 177
$ws.Range("K7").Value = 499

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 142
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 347

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 184
$ws.Range("K4").Value = 80
$ws.Range("K6").Value = 145
$ws.Range("K7").Value = 618
$ws.Range("K8").Value = 1389
$ws.Range("K11").Value = 395
$ws.Range("K15").Value = 217
$ws.Range("K18").Value = 139
$ws.Range("K19").Value = 614
$ws.Range("K20").Value = 501
$ws.Range("K23").Value = 216
$ws.Range("K24").Value = 62
$ws.Range("K25").Value = 101
$ws.Range("K26").Value = 29
$ws.Range("K27").Value = 197
$ws.Range("K29").Value = 1146
$ws.Range("K31").Value = 236
$ws.Range("K33").Value = 920
$ws.Range("K37").Value = 717
$ws.Range("K38").Value = 21
$ws.Range("K41").Value = 148
$ws.Range("K42").Value = 782
$ws.Range("K43").Value = 177
$ws.Range("K46").Value = 42
$ws.Range("K47").Value = 147
$ws.Range("K48").Value = 265
$ws.Range("I51").Value = 245
$ws.Range("K52").Value = 559
$ws.Range("K53").Value = 269
$ws.Range("K60").Value = 125
$ws.Range("K61").Value = 17
$ws.Range("C63").Value = 229
$ws.Range("G63").Value = 200
$ws.Range("J63").Value = 102
$ws.Range("K63").Value = 60
$ws.Range("K64").Value = 131
$ws.Range("K65").Value = 499
$ws.Range("K67").Value = 824
$ws.Range("K77").Value = 146
$ws.Range("K79").Value = 529
$ws.Range("K85").Value = 977
$ws.Range("K89").Value = 311
$ws.Range("K91").Value = 242
$ws.Range("K93").Value = 80
$ws.Range("K94").Value = 285
$ws.Range("K97").Value = 165
$ws.Range("K98").Value = 103
$ws.Range("K99").Value = 347
$ws.Range("C101").Value = 22626
$ws.Range("G101").Value = 20028
$ws.Range("I101").Value = 20782
$ws.Range("J101").Value = 23205
$ws.Range("K101").Value = 21108

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 75
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 228
$ws.Range("K3").Value = 296
$ws.Range("K7").Value = 824

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K6").Value = 321
$ws.Range("K7").Value = 1146

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 64
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 181
$ws.Range("K3").Value = 186
$ws.Range("K6").Value = 200
$ws.Range("K7").Value = 614

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 54
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K5").Value = 10
$ws.Range("K6").Value = 292
$ws.Range("K7").Value = 782

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 116
$ws.Range("K7").Value = 242

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 174
$ws.Range("K6").Value = 133
$ws.Range("K7").Value = 529

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 167
$ws.Range("K3").Value = 164
$ws.Range("K7").Value = 501

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 139

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 202
$ws.Range("K3").Value = 204
$ws.Range("K5").Value = 23
$ws.Range("K6").Value = 167
$ws.Range("K7").Value = 618

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 125
$ws.Range("K7").Value = 285

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 139
$ws.Range("K7").Value = 395

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 88
$ws.Range("K3").Value = 94
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 54
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I3").Value = 67
$ws.Range("I7").Value = 245

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 318
$ws.Range("K3").Value = 338
$ws.Range("K6").Value = 238
$ws.Range("K7").Value = 977

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 63
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 160
$ws.Range("K6").Value = 199
$ws.Range("K7").Value = 559

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 17
